$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new blank row at position 13 (shifts old rows 13-15 down to 14-16)
$ws.Rows.Item(13).Insert()

# Step 2: populate the new row 13 with the inserted match data (Poland - Division 1)
$ws.Cells.Item(13, 1).Value = "YP0aYXu2"
$ws.Cells.Item(13, 2).Value = "16/11/2024"
$ws.Cells.Item(13, 3).Value = "13:30"
$ws.Cells.Item(13, 4).Value = "POLAND - DIVISION 1"
$ws.Cells.Item(13, 5).Value = "Chrobry Glogow"
$ws.Cells.Item(13, 6).Value = "Wisla"
$ws.Cells.Item(13, 7).Value = 6
$ws.Cells.Item(13, 8).Value = 4.5
$ws.Cells.Item(13, 9).Value = 1.4
$ws.Cells.Item(13, 10).Value = 6.5
$ws.Cells.Item(13, 11).Value = 2.5
$ws.Cells.Item(13, 12).Value = 1.91
$ws.Cells.Item(13, 13).Value = 1.03
$ws.Cells.Item(13, 14).Value = 17
$ws.Cells.Item(13, 15).Value = 1.17
$ws.Cells.Item(13, 16).Value = 5
$ws.Cells.Item(13, 17).Value = 1.57
$ws.Cells.Item(13, 18).Value = 2.35
$ws.Cells.Item(13, 19).Value = 1.29
$ws.Cells.Item(13, 20).Value = 3.5
$ws.Cells.Item(13, 21).Value = 1.8
$ws.Cells.Item(13, 22).Value = 1.91
$ws.Cells.Item(13, 23).Value = 19
$ws.Cells.Item(13, 24).Value = 34
$ws.Cells.Item(13, 25).Value = 19
$ws.Cells.Item(13, 26).Value = 67
$ws.Cells.Item(13, 27).Value = 41
$ws.Cells.Item(13, 28).Value = 41
$ws.Cells.Item(13, 29).Value = 15
$ws.Cells.Item(13, 30).Value = 9.5
$ws.Cells.Item(13, 31).Value = 17
$ws.Cells.Item(13, 32).Value = 51
$ws.Cells.Item(13, 33).Value = 201
$ws.Cells.Item(13, 34).Value = 8.5
$ws.Cells.Item(13, 35).Value = 7.5
$ws.Cells.Item(13, 36).Value = 8.5
$ws.Cells.Item(13, 37).Value = 10
$ws.Cells.Item(13, 38).Value = 11
$ws.Cells.Item(13, 39).Value = 23
$ws.Cells.Item(13, 40).Value = 8.5
$ws.Cells.Item(13, 41).Value = 34
$ws.Cells.Item(13, 42).Value = 34
$ws.Cells.Item(13, 43).Value = 126
$ws.Cells.Item(13, 44).Value = 126
$ws.Cells.Item(13, 45).Value = 201
$ws.Cells.Item(13, 46).Value = 3.5
$ws.Cells.Item(13, 47).Value = 8.5
$ws.Cells.Item(13, 48).Value = 51
$ws.Cells.Item(13, 49).Value = 3.5
$ws.Cells.Item(13, 50).Value = 7
$ws.Cells.Item(13, 51).Value = 17
$ws.Cells.Item(13, 52).Value = 19
$ws.Cells.Item(13, 53).Value = 41
$ws.Cells.Item(13, 54).Value = 101
$ws.Cells.Item(13, 55).Value = 81
$ws.Cells.Item(13, 56).Value = 81

# Step 3: update odds values in row 6 (unrelated data refresh)
$ws.Cells.Item(6, 7).Value = 4.5
$ws.Cells.Item(6, 9).Value = 1.73
$ws.Cells.Item(6, 10).Value = 4.75
$ws.Cells.Item(6, 11).Value = 2.25
$ws.Cells.Item(6, 12).Value = 2.3
$ws.Cells.Item(6, 14).Value = 12
$ws.Cells.Item(6, 17).Value = 1.8
$ws.Cells.Item(6, 18).Value = 2
$ws.Cells.Item(6, 19).Value = 1.36
$ws.Cells.Item(6, 20).Value = 3
$ws.Cells.Item(6, 21).Value = 1.75
$ws.Cells.Item(6, 22).Value = 2
$ws.Cells.Item(6, 26).Value = 51
$ws.Cells.Item(6, 28).Value = 41
$ws.Cells.Item(6, 29).Value = 12
$ws.Cells.Item(6, 31).Value = 15
$ws.Cells.Item(6, 33).Value = 201
$ws.Cells.Item(6, 35).Value = 8.5
$ws.Cells.Item(6, 37).Value = 13
$ws.Cells.Item(6, 40).Value = 6.5
$ws.Cells.Item(6, 43).Value = 81
$ws.Cells.Item(6, 44).Value = 101
$ws.Cells.Item(6, 45).Value = 201
$ws.Cells.Item(6, 46).Value = 3
$ws.Cells.Item(6, 47).Value = 8
$ws.Cells.Item(6, 49).Value = 3.75
$ws.Cells.Item(6, 56).Value = 126

# Step 4: update odds values in row 10
$ws.Cells.Item(10, 15).Value = 1.29
$ws.Cells.Item(10, 16).Value = 3.5
$ws.Cells.Item(10, 17).Value = 1.93
$ws.Cells.Item(10, 18).Value = 1.93

# Step 5: update odds values in row 12
$ws.Cells.Item(12, 7).Value = 2.9
$ws.Cells.Item(12, 8).Value = 3.75
$ws.Cells.Item(12, 9).Value = 2.2
$ws.Cells.Item(12, 10).Value = 3.4
$ws.Cells.Item(12, 11).Value = 2.38
$ws.Cells.Item(12, 12).Value = 2.75
$ws.Cells.Item(12, 13).Value = 1.03
$ws.Cells.Item(12, 14).Value = 17
$ws.Cells.Item(12, 15).Value = 1.17
$ws.Cells.Item(12, 16).Value = 5
$ws.Cells.Item(12, 17).Value = 1.57
$ws.Cells.Item(12, 18).Value = 2.35
$ws.Cells.Item(12, 21).Value = 1.5
$ws.Cells.Item(12, 22).Value = 2.5
$ws.Cells.Item(12, 23).Value = 13
$ws.Cells.Item(12, 27).Value = 21
$ws.Cells.Item(12, 28).Value = 23
$ws.Cells.Item(12, 29).Value = 17
$ws.Cells.Item(12, 31).Value = 12
$ws.Cells.Item(12, 34).Value = 11
$ws.Cells.Item(12, 35).Value = 13
$ws.Cells.Item(12, 37).Value = 21
$ws.Cells.Item(12, 40).Value = 5.5
$ws.Cells.Item(12, 42).Value = 21
$ws.Cells.Item(12, 49).Value = 4.5
$ws.Cells.Item(12, 55).Value = 351
